{"js": "// Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" /\n// \"\u00a9 2020 . Contact: ...\" footer block (and the blank paragraph that\n// separates it from the bibliography entry above), as built by the\n// 2020-10-28 site regeneration.\n\nconst body = context.document.body;\n\n// Locate the unique copyright paragraph - this is the most distinctive\n// anchor text in the block being removed.\nconst copyrightResults = body.search(\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\",\n  { matchCase: true }\n);\ncopyrightResults.load(\"items\");\nawait context.sync();\n\nif (copyrightResults.items.length === 0) {\n  throw new Error(\"Could not find the copyright footer paragraph to remove.\");\n}\n\nconst copyrightRange = copyrightResults.items[0];\nconst copyrightParagraphs = copyrightRange.paragraphs;\ncopyrightParagraphs.load(\"items\");\nawait context.sync();\n\nconst copyrightPara = copyrightParagraphs.items[0];\nconst jupiterPara = copyrightPara.getPrevious();\nconst blankPara = jupiterPara.getPrevious();\n\njupiterPara.load(\"text\");\nblankPara.load(\"text\");\nawait context.sync();\n\n// Sanity-check the chain before deleting anything, so we never remove the\n// wrong paragraphs if the document shape differs from what we expect.\nif (jupiterPara.text !== \"Ver no Jupiter Salvar em pdf Salvar em docx\") {\n  throw new Error(\"Unexpected paragraph before the copyright line: \" + jupiterPara.text);\n}\nif (blankPara.text !== \"\") {\n  throw new Error(\"Unexpected paragraph before the 'Ver no Jupiter' line: \" + blankPara.text);\n}\n\n// Delete in reverse document order so earlier deletes never invalidate the\n// object identities of paragraphs we still need to remove.\ncopyrightPara.delete();\njupiterPara.delete();\nblankPara.delete();\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" /\n# \"\u00a9 2020 . Contact: ...\" footer block (and the blank paragraph that\n# separates it from the bibliography entry above), as built by the\n# 2020-10-28 site regeneration.\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphIndexByText($doc, $text) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        $p = $doc.Paragraphs.Item($i)\n        # Paragraph.Range.Text carries a trailing paragraph mark (CR, or\n        # cell-mark 0x07 inside a table); strip it before comparing.\n        $ptext = $p.Range.Text.TrimEnd([char]13, [char]7)\n        if ($ptext -eq $text) {\n            return $i\n        }\n    }\n    return -1\n}\n\n$jupiterIdx = Find-ParagraphIndexByText $d \"Ver no Jupiter Salvar em pdf Salvar em docx\"\nif ($jupiterIdx -eq -1) {\n    throw \"Could not find the 'Ver no Jupiter ...' paragraph to remove.\"\n}\n\n$blankIdx = $jupiterIdx - 1\n$copyrightIdx = $jupiterIdx + 1\n\n# Sanity-check the surrounding paragraphs before deleting anything, so we\n# never remove the wrong block if the document shape differs from expected.\n$blankText = $d.Paragraphs.Item($blankIdx).Range.Text.TrimEnd([char]13, [char]7)\nif ($blankText -ne \"\") {\n    throw \"Unexpected paragraph before the 'Ver no Jupiter' line: $blankText\"\n}\n\n$copyrightText = $d.Paragraphs.Item($copyrightIdx).Range.Text.TrimEnd([char]13, [char]7)\nif ($copyrightText -notlike \"*Contact: luizeleno@usp.br*\") {\n    throw \"Unexpected paragraph after the 'Ver no Jupiter' line: $copyrightText\"\n}\n\n# Delete highest index first so the lower indices stay valid while we work.\n$d.Paragraphs.Item($copyrightIdx).Range.Delete()\n$d.Paragraphs.Item($jupiterIdx).Range.Delete()\n$d.Paragraphs.Item($blankIdx).Range.Delete()\n"}
